$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting existing rows 46-55 down to 47-56
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new weekly record
$ws.Cells.Item(46, 1).Value = 10
$ws.Cells.Item(46, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(46, 3).Value = "La Araucanía"
$ws.Cells.Item(46, 4).Value = 44627
$ws.Cells.Item(46, 4).NumberFormat = $ws.Cells.Item(47, 4).NumberFormat
$ws.Cells.Item(46, 5).Value = 9
$ws.Cells.Item(46, 6).Value = 100114002
$ws.Cells.Item(46, 7).Value = "Camote"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 20
$ws.Cells.Item(46, 11).Value = 18000
$ws.Cells.Item(46, 12).Value = 18000
$ws.Cells.Item(46, 13).Value = 18000
$ws.Cells.Item(46, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(46, 15).Value = "Perú"
$ws.Cells.Item(46, 16).Value = 900
$ws.Cells.Item(46, 17).Value = 20
$ws.Cells.Item(46, 18).Value = "Hortaliza"
